$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 carries the "excel formating" task header; row 29 (sub="c", main count=1)
# gets fleshed out with the new sub-task description, wrapped over two lines,
# plus start/finish dates of 2020-07-08.
$ws.Range("D29").Value = "upper part table(General trade figures, % change, rank)"
$ws.Range("D29").WrapText = $true
$ws.Range("D29").HorizontalAlignment = -4131  # xlLeft

$ws.Range("E29").Value = 20200708
$ws.Range("F29").Value = 20200708

# Two-line text needs the taller row height used by the other wrapped rows.
$ws.Rows("29").RowHeight = 28.8

# Move the active selection to where the author left off editing.
$ws.Range("E30").Select()
